$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.291.96"
$ws.Range("E2").Value = "  +0.20%  "

# Row 3
$ws.Range("D3").Value = "3.328.48"
$ws.Range("E3").Value = "  -0.24%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.11%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.54%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").Value = "3.318.93"
$ws.Range("E9").Value = "  -0.32%  "

# Row 10
$ws.Range("E10").Value = "  +7.32%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.634"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.76%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.09"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.62%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000280"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.42%  "

# Row 14
$ws.Range("E14").Value = "  +0.70%  "

# Row 15
$ws.Range("D15").Value = "3.859.51"
$ws.Range("E15").Value = "  -0.07%  "

# Row 16
$ws.Range("E16").Value = "  +2.90%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.11"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.31%  "

# Row 18
$ws.Range("D18").Value = "3.321.45"
$ws.Range("E18").Value = "  -0.20%  "

# Row 19
$ws.Range("D19").Value = "64.338.14"
$ws.Range("E19").Value = "  +0.49%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.71"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.986"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "452.44"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.43%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.78%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.06"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.63"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.43%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.87"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.39%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.86"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.78%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.53"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.11%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "31.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.38%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.59"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.33%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.50"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.64%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.39"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "62.35"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "571.86"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.84%  "

# Row 35
$ws.Range("E35").Value = "  -0.17%  "

# Row 37
$ws.Range("E37").Value = "  +0.48%  "

# Row 38
$ws.Range("E38").Value = "  -0.37%  "

# Row 39
$ws.Range("E39").Value = "  -0.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.366"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.29%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0729"
$ws.Range("E41").Value = "  -3.21%  "

# Row 42
$ws.Range("D42").Value = "3.071.83"
$ws.Range("E42").Value = "  -0.62%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0413"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.05%  "

# Row 44
$ws.Range("E44").Value = "  -1.93%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.18"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.00%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.134"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.65%  "

# Row 47
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.44"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.22%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.18%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.23"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.84%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.16"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.42%  "
